$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated vm_pu values for the 380 kV case (rows 2-25, columns B-F and I-M).
$updates = @{
    "B2" = 1.02;
    "C2" = 1.044347406509429;
    "D2" = 1.051968282271224;
    "E2" = 1.052089468472426;
    "F2" = 1.062542699770794;
    "I2" = 1.042932735327767;
    "J2" = 1.049412977925657;
    "K2" = 1.054718315282301;
    "L2" = 1.054839165913459;
    "M2" = 1.065263763622047;
    "B3" = 1.02;
    "C3" = 1.045343851258234;
    "D3" = 1.052749651937033;
    "E3" = 1.052953231347063;
    "F3" = 1.063439417389199;
    "I3" = 1.043152968992481;
    "J3" = 1.050056408429247;
    "K3" = 1.055312517614185;
    "L3" = 1.055515574008128;
    "M3" = 1.065975111797181;
    "B4" = 1.02;
    "C4" = 1.045988991601382;
    "D4" = 1.053255244276344;
    "E4" = 1.053512762115651;
    "F4" = 1.064020079858497;
    "I4" = 1.043293857447577;
    "J4" = 1.050472513907413;
    "K4" = 1.055696341392797;
    "L4" = 1.055953230913457;
    "M4" = 1.06643519045166;
    "B5" = 1.02;
    "C5" = 1.046260297376032;
    "D5" = 1.053467792565738;
    "E5" = 1.053748135687779;
    "F5" = 1.064264291108597;
    "I5" = 1.043352699109338;
    "J5" = 1.050647387016301;
    "K5" = 1.055857540747686;
    "L5" = 1.056137215074529;
    "M5" = 1.066628555876497;
    "B6" = 1.02;
    "C6" = 1.046305855993346;
    "D6" = 1.05350348016159;
    "E6" = 1.053787664521422;
    "F6" = 1.064305301122513;
    "I6" = 1.043362556121699;
    "J6" = 1.050676745575575;
    "K6" = 1.055884597407832;
    "L6" = 1.056168106397653;
    "M6" = 1.066661019722324;
    "B7" = 1.02;
    "C7" = 1.045992616454539;
    "D7" = 1.053258084369326;
    "E7" = 1.053515906613406;
    "F7" = 1.064023342626926;
    "I7" = 1.043294645217493;
    "J7" = 1.050474850798021;
    "K7" = 1.055698495977046;
    "L7" = 1.055955689346773;
    "M7" = 1.066437774414216;
    "B8" = 1.02;
    "C8" = 1.044684081793975;
    "D8" = 1.05223235030321;
    "E8" = 1.052381252399294;
    "F8" = 1.062845660666198;
    "I8" = 1.043007498856112;
    "J8" = 1.049630476840282;
    "K8" = 1.054919265811351;
    "L8" = 1.055067765658036;
    "M8" = 1.065504210143885;
    "B9" = 1.02;
    "C9" = 1.042381173604827;
    "D9" = 1.050424885170535;
    "E9" = 1.050386636631545;
    "F9" = 1.060773755391653;
    "I9" = 1.042489149843759;
    "J9" = 1.048140804726659;
    "K9" = 1.053541114161126;
    "L9" = 1.053502986677744;
    "M9" = 1.063857580832341;
    "B10" = 1.02;
    "C10" = 1.040847891787882;
    "D10" = 1.049219988837329;
    "E10" = 1.049060184164096;
    "F10" = 1.059394792829574;
    "I10" = 1.042135309414604;
    "J10" = 1.047146545637738;
    "K10" = 1.052619009500177;
    "L10" = 1.052459760146617;
    "M10" = 1.062758835441816;
    "B11" = 1.02;
    "C11" = 1.040184444223569;
    "D11" = 1.048698289728792;
    "E11" = 1.048486610842391;
    "F11" = 1.058798249518441;
    "I11" = 1.041980137350508;
    "J11" = 1.046715759332396;
    "K11" = 1.052218948258388;
    "L11" = 1.052008033421842;
    "M11" = 1.062282844009538;
    "B12" = 1.02;
    "C12" = 1.039938082001035;
    "D12" = 1.048504512981638;
    "E12" = 1.048273679897675;
    "F12" = 1.058576751582396;
    "I12" = 1.041922206088461;
    "J12" = 1.046555706998801;
    "K12" = 1.052070230828977;
    "L12" = 1.051840242471266;
    "M12" = 1.062106006415741;
    "B13" = 1.02;
    "C13" = 1.039990924278012;
    "D13" = 1.048546078482521;
    "E13" = 1.048319348872116;
    "F13" = 1.058624259784606;
    "I13" = 1.041934645811293;
    "J13" = 1.046590040525406;
    "K13" = 1.052102136496815;
    "L13" = 1.051876234164422;
    "M13" = 1.062143940165664;
    "B14" = 1.02;
    "C14" = 1.040164078378432;
    "D14" = 1.048682271962793;
    "E14" = 1.048469007455163;
    "F14" = 1.058779938692956;
    "I14" = 1.041975354717443;
    "J14" = 1.046702530150093;
    "K14" = 1.052206657608276;
    "L14" = 1.051994163754882;
    "M14" = 1.062268227230085;
    "B15" = 1.02;
    "C15" = 1.040270773858907;
    "D15" = 1.048766186029015;
    "E15" = 1.048561232946199;
    "F15" = 1.058875868880287;
    "I15" = 1.042000397949517;
    "J15" = 1.046771833560096;
    "K15" = 1.052271041054696;
    "L15" = 1.052066824156111;
    "M15" = 1.062344800224944;
    "B16" = 1.02;
    "C16" = 1.040891932646056;
    "D16" = 1.049254612999868;
    "E16" = 1.049098267064021;
    "F16" = 1.059434395272673;
    "I16" = 1.04214556649175;
    "J16" = 1.047175129977499;
    "K16" = 1.052645543808796;
    "L16" = 1.05248973977892;
    "M16" = 1.062790420718126;
    "B17" = 1.02;
    "C17" = 1.041281696351276;
    "D17" = 1.04956099892114;
    "E17" = 1.049435346510898;
    "F17" = 1.059784893976193;
    "I17" = 1.042236103369134;
    "J17" = 1.047428036617396;
    "K17" = 1.052880250266552;
    "L17" = 1.052755023461918;
    "M17" = 1.063069886488515;
    "B18" = 1.02;
    "C18" = 1.041509084580057;
    "D18" = 1.049739711360174;
    "E18" = 1.049632035334215;
    "F18" = 1.059989387561974;
    "I18" = 1.042288723156613;
    "J18" = 1.047575527106933;
    "K18" = 1.053017074883754;
    "L18" = 1.052909758621361;
    "M18" = 1.063232872134748;
    "B19" = 1.02;
    "C19" = 1.041586625815506;
    "D19" = 1.049800648108051;
    "E19" = 1.04969911401084;
    "F19" = 1.060059123649071;
    "I19" = 1.042306633097041;
    "J19" = 1.047625813173805;
    "K19" = 1.053063715677081;
    "L19" = 1.052962519230107;
    "M19" = 1.063288442264966;
    "B20" = 1.02;
    "C20" = 1.041239873667849;
    "D20" = 1.049528126322014;
    "E20" = 1.049399173198278;
    "F20" = 1.059747283227361;
    "I20" = 1.04222640914782;
    "J20" = 1.047400904778402;
    "K20" = 1.052855076307222;
    "L20" = 1.052726561060513;
    "M20" = 1.063039904725557;
    "B21" = 1.02;
    "C21" = 1.040113086826385;
    "D21" = 1.048642166215259;
    "E21" = 1.048424933418494;
    "F21" = 1.058734092783895;
    "I21" = 1.041963375055491;
    "J21" = 1.046669405832986;
    "K21" = 1.05217588196515;
    "L21" = 1.051959436409813;
    "M21" = 1.062231628695474;
    "B22" = 1.02;
    "C22" = 1.039405046594322;
    "D22" = 1.048085160511318;
    "E22" = 1.04781308345524;
    "F22" = 1.058097551015612;
    "I22" = 1.041796297547179;
    "J22" = 1.046209257466639;
    "K22" = 1.051748169782242;
    "L22" = 1.051477117838632;
    "M22" = 1.061723241978731;
    "B23" = 1.02;
    "C23" = 1.039780352546444;
    "D23" = 1.048380436292671;
    "E23" = 1.048137370606701;
    "F23" = 1.058434946886321;
    "I23" = 1.0418850292216;
    "J23" = 1.046453211979522;
    "K23" = 1.051974971858643;
    "L23" = 1.051732803447034;
    "M23" = 1.061992765219269;
    "B24" = 1.02;
    "C24" = 1.041258771397168;
    "D24" = 1.049542980026275;
    "E24" = 1.049415518126185;
    "F24" = 1.059764277740225;
    "I24" = 1.042230790132662;
    "J24" = 1.047413164567934;
    "K24" = 1.052866451569758;
    "L24" = 1.052739421995615;
    "M24" = 1.063053452261754;
    "B25" = 1.02;
    "C25" = 1.042976182881879;
    "D25" = 1.05089214946898;
    "E25" = 1.05090171767139;
    "F25" = 1.061308991305862;
    "I25" = 1.042624616542084;
    "J25" = 1.048526126073632;
    "K25" = 1.053897991741942;
    "L25" = 1.053907530809498;
    "M25" = 1.06428345351231
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output "Updated $($updates.Count) cells"
